# Edit script: applies "add func HTofRS, SM_CS_ExLos; add data d_rs; add if calculation H0_rs"
# to the turbine preliminary-calculation worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the turbine nameplate / input data block (row 3)
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "Т-100-130"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 550
$ws.Range("G3").Value = 3.6
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 180
$ws.Range("K3").Value = 130

# ---------------------------------------------------------------------------
# 2. Recalculated steam-table results for the existing blocks (1.1 - 3)
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = 3470.996181930883
$ws.Range("E9").Value = 6.608253838576329
$ws.Range("F9").Value = 0.02685306336725179

$ws.Range("C14").Value = 341.069773726435
$ws.Range("D14").Value = 3080.291195557765
$ws.Range("F14").Value = 0.07308081411324667

$ws.Range("C18").Value = 390.7049863731177

$ws.Range("C22").Value = 230.7382694367192

$ws.Range("H24").Value = 43172.9009942295

$ws.Range("F27").Value = 12.35

$ws.Range("C32").Value = 547.4190007152333
$ws.Range("D32").Value = 3470.996181930883
$ws.Range("E32").Value = 6.630048689665621
$ws.Range("F32").Value = 0.02827546851902559

# ---------------------------------------------------------------------------
# 3. New section 4 "Параметры пара на выходе из регулирующей ступени"
# ---------------------------------------------------------------------------

# 4. title
$ws.Range("B26").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Range("B34").Value = "4. Параметры пара на выходе из регулирующей ступени"

# 4.1 title
$ws.Range("B26").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Range("B35").Value = "4.1 Параметры пара после регулирующей ступени (без учета потерь)"

# 4.1.1 subsection header (style like B11 / B16 / B20)
$ws.Range("B20").Copy()
$ws.Range("B36").PasteSpecial(-4122)
$ws.Range("B36").Value = "4.1.1 Располагаемый теплоперепад рег.ступени"

# row 37: u/cf ratio label + value
$ws.Rows.Item(37).RowHeight = 18
$ws.Range("B26").Copy()
$ws.Range("B37").PasteSpecial(-4122)
$ws.Range("B37").Value = "Отношение скоростей, u/cф="

$ws.Range("H21").Copy()
$ws.Range("E37").PasteSpecial(-4122)
$ws.Range("E37").Value = 0.24

# row 38: H0рс label (new "left aligned" style) + value
$ws.Rows.Item(38).RowHeight = 18.75
$ws.Range("D38").HorizontalAlignment = -4131
$ws.Range("D38").Value = "H0рс, кДж/кг"

$ws.Range("E38").Value = 180

# row 41: table header (Р, МПа / t, oC / h, кДж/кг / S, кДж/кг / v, м3/кг)
$ws.Rows.Item(41).RowHeight = 17.25
$ws.Range("A31").Copy()
$ws.Range("A41").PasteSpecial(-4122)

$ws.Range("B31").Copy()
$ws.Range("B41").PasteSpecial(-4122)
$ws.Range("B41").Value = "Р, МПа"

$ws.Range("C31").Copy()
$ws.Range("C41").PasteSpecial(-4122)
$ws.Range("C41").Value = "t, oC"

$ws.Range("D31").Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("D41").Value = "h, кДж/кг"

$ws.Range("E31").Copy()
$ws.Range("E41").PasteSpecial(-4122)
$ws.Range("E41").Value = "S, кДж/кг"

$ws.Range("F31").Copy()
$ws.Range("F41").PasteSpecial(-4122)
$ws.Range("F41").Value = "v, м3/кг"

# row 42: "Индекс 2рсt" data row
$ws.Range("A32").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("A42").Value = "Индекс 2рсt"

$ws.Range("B32").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$ws.Range("B42").Value = 7.148328643637858

$ws.Range("C32").Copy()
$ws.Range("C42").PasteSpecial(-4122)
$ws.Range("C42").Value = 452.1459474068108

$ws.Range("D32").Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("D42").Value = 3290.996181930883

$ws.Range("E32").Copy()
$ws.Range("E42").PasteSpecial(-4122)
$ws.Range("E42").Value = 6.630048689665621

$ws.Range("F32").Copy()
$ws.Range("F42").PasteSpecial(-4122)
$ws.Range("F42").Value = 0.0433555151468901

# ---------------------------------------------------------------------------
# 4. Restore view / selection state
# ---------------------------------------------------------------------------
$ws.Range("K41:L41").Select()
$excel.ActiveWindow.ScrollRow = 28
